# Fixed update to excel issue
#
# 1. Rename "Requested quantity" header to "Weekly_PO_Qty" on the
#    "Weekly Quantity" sheet and to "Monthly_PO_Qty" on the "Monthly Trend"
#    sheet.
# 2. Add a new "PO Forecast" worksheet (after the existing sheets) that
#    contains the ds / PO_Forecast / yhat_lower / yhat_upper forecast table.

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet at the end of the workbook ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows
$data = @()
$data += ,@(45347.99999999999, 153, 50.69909505993647, 251.0156420253972)
$data += ,@(45361.99999999999, 146, 44.4866965335159, 245.1896174558467)
$data += ,@(45424.99999999999, 115, 15.0176257940486, 208.2020863486384)
$data += ,@(45431.99999999999, 111, 12.07627030762826, 220.6182299572116)
$data += ,@(45445.99999999999, 104, -1.910908150980864, 203.7220365366918)
$data += ,@(45459.99999999999, 97, -2.599708534968669, 202.5102576064041)
$data += ,@(45522.99999999999, 66, -40.87810640052459, 165.4639542811665)
$data += ,@(45529.99999999999, 63, -34.49028562732039, 175.7857606484004)
$data += ,@(45536.99999999999, 59, -47.25446142215318, 167.7465784967583)
$data += ,@(45543.99999999999, 56, -43.42746904944066, 156.9528344323664)
$data += ,@(45564.99999999999, 45, -53.93469500794424, 141.5442567440089)
$data += ,@(45578.99999999999, 38, -60.81936914884518, 137.0317911432386)
$data += ,@(45585.99999999999, 35, -64.07505619933787, 131.2743028662732)
$data += ,@(45592.99999999999, 31, -72.86258484537909, 138.1838551785591)
$data += ,@(45599.99999999999, 28, -71.10832532723147, 131.6756231281621)
$data += ,@(45606.99999999999, 24, -80.21597195762203, 119.8820176214026)
$data += ,@(45613.99999999999, 21, -84.77193095553113, 117.6596410021138)
$data += ,@(45620.99999999999, 17, -83.51218365319836, 123.1686657721862)
$data += ,@(45627.99999999999, 14, -83.51651717308471, 117.4474478163459)
$data += ,@(45634.99999999999, 10, -96.28793593319349, 100.6107935160628)
$data += ,@(45641.99999999999, 7, -101.3465757649452, 103.9178532216408)

$r = 2
foreach ($row in $data) {
    $dsCell = $wsForecast.Cells.Item($r, 1)
    $dsCell.Value = $row[0]
    $dsCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]

    $r++
}
